$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Reporte generado" timestamp in A1
$ws.Range("A1").Value = "Reporte generado: 2019-10-12 08:56:57"

# Add a new data row (row 5) with a thin border around it
$rowRange = $ws.Range("A5:E5")
$rowRange.Borders.Color = 0
$rowRange.Borders.LineStyle = 1

$ws.Range("A5").Value = "CARTAGO"
$ws.Range("B5").Value = "SEDE1"
$ws.Range("C5").Value = "JESUS CAMARA"
$ws.Range("D5").Value = "2019-10-12 08:56:51"
$ws.Range("E5").Value = "CONECTADO"

# New green fill for the "CONECTADO" status cell
$ws.Range("E5").Interior.Color = 3580485

# Adjust column widths to fit the new, wider content
$ws.Columns.Item(1).ColumnWidth = 8.5
$ws.Columns.Item(2).ColumnWidth = 6.166667
$ws.Columns.Item(3).ColumnWidth = 14.5
$ws.Columns.Item(4).ColumnWidth = 22.666667
$ws.Columns.Item(5).ColumnWidth = 10.833333

# Update the selected range to include the new row
$ws.Range("A4:E5").Select()
